$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 53) entirely - the series now ends one period earlier
$ws.Rows.Item(53).Delete()

# Bulk-write the corrected forecast data (dates shifted + recalculated AR(2) forecasts) into A2:E52
$data = New-Object "object[,]" 51,5
$data[0,0] = 39583
$data[0,1] = 2008
$data[0,2] = $null
$data[0,3] = 2009
$data[0,4] = $null
$data[1,0] = 39765
$data[1,1] = 2008
$data[1,2] = $null
$data[1,3] = 2009
$data[1,4] = $null
$data[2,0] = 39948
$data[2,1] = 2009
$data[2,2] = $null
$data[2,3] = 2010
$data[2,4] = $null
$data[3,0] = 40130
$data[3,1] = 2009
$data[3,2] = -4.715480642250625
$data[3,3] = 2010
$data[3,4] = $null
$data[4,0] = 40310
$data[4,1] = 2010
$data[4,2] = $null
$data[4,3] = 2011
$data[4,4] = $null
$data[5,0] = 40494
$data[5,1] = 2010
$data[5,2] = 6.130685532900881
$data[5,3] = 2011
$data[5,4] = $null
$data[6,0] = 40676
$data[6,1] = 2011
$data[6,2] = $null
$data[6,3] = 2012
$data[6,4] = $null
$data[7,0] = 40862
$data[7,1] = 2011
$data[7,2] = 8.703939237319025
$data[7,3] = 2012
$data[7,4] = $null
$data[8,0] = 41044
$data[8,1] = 2012
$data[8,2] = $null
$data[8,3] = 2013
$data[8,4] = $null
$data[9,0] = 41228
$data[9,1] = 2012
$data[9,2] = 2.688274587589135
$data[9,3] = 2013
$data[9,4] = 1.490702606731831
$data[10,0] = 41409
$data[10,1] = 2013
$data[10,2] = 0.4454453461194552
$data[10,3] = 2014
$data[10,4] = 2.887668087172179
$data[11,0] = 41592
$data[11,1] = 2013
$data[11,2] = 0.9946838291217786
$data[11,3] = 2014
$data[11,4] = 3.498411812952029
$data[12,0] = 41774
$data[12,1] = 2014
$data[12,2] = 5.461771395837989
$data[12,3] = 2015
$data[12,4] = 3.667647087004666
$data[13,0] = 41957
$data[13,1] = 2014
$data[13,2] = 5.562499360312567
$data[13,3] = 2015
$data[13,4] = 4.543069198269034
$data[14,0] = 42137
$data[14,1] = 2015
$data[14,2] = 3.857137494611718
$data[14,3] = 2016
$data[14,4] = 3.60208822706134
$data[15,0] = 42321
$data[15,1] = 2015
$data[15,2] = 4.195080504802551
$data[15,3] = 2016
$data[15,4] = 4.42512316868644
$data[16,0] = 42503
$data[16,1] = 2016
$data[16,2] = 4.325618632128836
$data[16,3] = 2017
$data[16,4] = 3.998755562728684
$data[17,0] = 42689
$data[17,1] = 2016
$data[17,2] = 4.230623896992025
$data[17,3] = 2017
$data[17,4] = 4.210645455310114
$data[18,0] = 42867
$data[18,1] = 2017
$data[18,2] = 4.439146757103352
$data[18,3] = 2018
$data[18,4] = 3.861679870292711
$data[19,0] = 43053
$data[19,1] = 2017
$data[19,2] = 4.933871867981643
$data[19,3] = 2018
$data[19,4] = 4.757571096183799
$data[20,0] = 43145
$data[20,1] = 2018
$data[20,2] = 5.787381971961936
$data[20,3] = 2019
$data[20,4] = 4.118094101621717
$data[21,0] = 43235
$data[21,1] = 2018
$data[21,2] = 5.723509166364238
$data[21,3] = 2019
$data[21,4] = 4.058053416301188
$data[22,0] = 43326
$data[22,1] = 2018
$data[22,2] = 5.222860865675738
$data[22,3] = 2019
$data[22,4] = 2.939060555390971
$data[23,0] = 43418
$data[23,1] = 2018
$data[23,2] = 5.456119081407906
$data[23,3] = 2019
$data[23,4] = 4.569144243718659
$data[24,0] = 43510
$data[24,1] = 2019
$data[24,2] = 3.466212706516147
$data[24,3] = 2020
$data[24,4] = 3.597750881470851
$data[25,0] = 43600
$data[25,1] = 2019
$data[25,2] = 4.674926984813466
$data[25,3] = 2020
$data[25,4] = 4.726969153629335
$data[26,0] = 43691
$data[26,1] = 2019
$data[26,2] = 3.038115835571786
$data[26,3] = 2020
$data[26,4] = 1.154972712087221
$data[27,0] = 43783
$data[27,1] = 2019
$data[27,2] = 3.346849276607955
$data[27,3] = 2020
$data[27,4] = 3.124801698476176
$data[28,0] = 43875
$data[28,1] = 2020
$data[28,2] = 3.370990011762443
$data[28,3] = 2021
$data[28,4] = 3.843786543692795
$data[29,0] = 43966
$data[29,1] = 2020
$data[29,2] = 0.4167846160013644
$data[29,3] = 2021
$data[29,4] = 1.125570778878981
$data[30,0] = 44068
$data[30,1] = 2020
$data[30,2] = -9.2489161297999
$data[30,3] = 2021
$data[30,4] = -10.65745199005891
$data[31,0] = 44159
$data[31,1] = 2020
$data[31,2] = -9.2489161297999
$data[31,3] = 2021
$data[31,4] = -4.101394328717845
$data[32,0] = 44251
$data[32,1] = 2021
$data[32,2] = -1.861534891151506
$data[32,3] = 2022
$data[32,4] = 3.108682697521514
$data[33,0] = 44341
$data[33,1] = 2021
$data[33,2] = -1.488064879190421
$data[33,3] = 2022
$data[33,4] = 3.63609986063671
$data[34,0] = 44432
$data[34,1] = 2021
$data[34,2] = -1.287084480507283
$data[34,3] = 2022
$data[34,4] = 4.124305474197043
$data[35,0] = 44525
$data[35,1] = 2021
$data[35,2] = -1.287084480507283
$data[35,3] = 2022
$data[35,4] = 3.801772939051373
$data[36,0] = 44617
$data[36,1] = 2022
$data[36,2] = 1.920033066224791
$data[36,3] = 2023
$data[36,4] = 2.774323849124349
$data[37,0] = 44706
$data[37,1] = 2022
$data[37,2] = 1.712986619197032
$data[37,3] = 2023
$data[37,4] = 2.431967849366434
$data[38,0] = 44798
$data[38,1] = 2022
$data[38,2] = 1.494343500592232
$data[38,3] = 2023
$data[38,4] = 1.882230726672129
$data[39,0] = 44890
$data[39,1] = 2022
$data[39,2] = 1.494343500592232
$data[39,3] = 2023
$data[39,4] = 0.9099262091262217
$data[40,0] = 44981
$data[40,1] = 2023
$data[40,2] = -0.07673633990846751
$data[40,3] = 2024
$data[40,4] = 2.180844122535164
$data[41,0] = 45071
$data[41,1] = 2023
$data[41,2] = -0.5717743519535134
$data[41,3] = 2024
$data[41,4] = 1.732880403074311
$data[42,0] = 45163
$data[42,1] = 2023
$data[42,2] = -0.6982718287330991
$data[42,3] = 2024
$data[42,4] = 1.415512869596025
$data[43,0] = 45254
$data[43,1] = 2023
$data[43,2] = -0.6982718287330991
$data[43,3] = 2024
$data[43,4] = 0.1232424362653362
$data[44,0] = 45345
$data[44,1] = 2024
$data[44,2] = -0.5700058398449448
$data[44,3] = 2025
$data[44,4] = 1.804123797928292
$data[45,0] = 45436
$data[45,1] = 2024
$data[45,2] = -0.2867681914691111
$data[45,3] = 2025
$data[45,4] = 2.013081730696564
$data[46,0] = 45534
$data[46,1] = 2024
$data[46,2] = -0.4137309550271362
$data[46,3] = 2025
$data[46,4] = 1.790585695398428
$data[47,0] = 45618
$data[47,1] = 2024
$data[47,2] = -0.4137309550271362
$data[47,3] = 2025
$data[47,4] = 1.743978804508384
$data[48,0] = 45713
$data[48,1] = 2025
$data[48,2] = 0.946259771301472
$data[48,3] = 2026
$data[48,4] = 1.770613036357038
$data[49,0] = 45800
$data[49,1] = 2025
$data[49,2] = 0.4432539413513181
$data[49,3] = 2026
$data[49,4] = 1.364481450639365
$data[50,0] = 45891
$data[50,1] = 2025
$data[50,2] = 0.2267356977060819
$data[50,3] = 2026
$data[50,4] = 0.9098136509666066

$ws.Range("A2:E52").Value2 = $data

Write-Output "Applied forecast data update"